# Populate the "PushNotificationInformations" sheet with the supermarket
# category test data (row 1 already existed; rows 2 and 3 are new).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Vegetables"
$ws.Range("B1").Value = "Fresh Vegetables"
$ws.Range("A2").Value = "Fruits"
$ws.Range("B2").Value = "Fresh Fruits"
$ws.Range("A3").Value = "Bakery"
$ws.Range("B3").Value = "Homemade items"

# Widen the columns to fit the new content (matches the saved widths of
# 13.85546875 / 17 characters as closely as the column-width quantization
# allows).
$ws.Columns.Item(1).ColumnWidth = 13.0
$ws.Columns.Item(2).ColumnWidth = 16.166666666666668

# Leave the selection where the author finished editing.
$ws.Range("B3").Select()
